$p = $ppt.ActivePresentation

# --- Change 1: reposition/resize "TextBox 103" (shape id 104, text "C") on slide 1 ---
$s1 = $p.Slides.Item(1)
$shpC = $s1.Shapes.Item(45)
$shpC.Left = 408.38807055856296
$shpC.Top = 252.91874015748033
$shpC.Width = 27.93748031496063
$shpC.Height = 29.081259842519685

# --- Change 2: merge the two text runs "1- " and "P_2" into a single run on slide 4 ---
$s4 = $p.Slides.Item(4)
$shpP2 = $s4.Shapes.Item(48)
# The concatenated text already reads as "1- P_2", and the engine preserves
# existing run boundaries when the new text is just a prefix-extension/
# truncation of the current text. Route through an unrelated placeholder
# value first so the final assignment forces a genuine single-run rebuild.
$shpP2.TextFrame.TextRange.Text = "ZZZZZZZZZZ"
$shpP2.TextFrame.TextRange.Text = "1- P_2"
